$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(0,1,2,3,4,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29,30,31,32,35,36,37,38,39,41,42,44,45,47,48,49,50,52,53,54,56,57,58,61,62,63,66,67,68,69,70,71,73,74,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
